$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 647: add SITE value that was missing ---
$ws.Range("A647").Value = "UIC"

# --- New rows 648-678 ---
$data = @(
    @{ Row=648; A="UIC"; B="081514530";  C=44181; D="PAF";    E="NICM";        G="TEE";    H="bedside" },
    @{ Row=649; A="UIC"; B="081514530";  C=44181; D="PAF";    E="NICM";        G="TEE" },
    @{ Row=650; A="UIC"; B="081550797";  C=44181; D="VSD";    E="CHD";         G="TTE" },
    @{ Row=651; A="UIC"; B="076462787";  C=44180;                              G="TTE" },
    @{ Row=652; A="UIC"; B="087054912";  C=44180;                              G="TTE" },
    @{ Row=653; A="UIC"; B="071434237";  C=44180;                              G="TTE" },
    @{ Row=654; A="UIC"; B="071991053";  C=44180;                              G="TTE" },
    @{ Row=655; A="UIC"; B="080911553";  C=44180; D="ASD";                     G="TTE" },
    @{ Row=656; A="UIC"; B="081550797";  C=44181;                              G="TTE" },
    @{ Row=657; A="UIC"; B="2000103647"; C=44181;                              G="TTE" },
    @{ Row=658; A="UIC"; B="080684682";  C=44181;                              G="TTE" },
    @{ Row=659; A="UIC"; B="200206802";  C=44181;                              G="TTE" },
    @{ Row=660; A="UIC"; B="051365492";  C=44181;                              G="TTE" },
    @{ Row=661; A="UIC"; B="081068545";  C=44181;                              G="TTE" },
    @{ Row=662; A="UIC"; B="007737547";  C=44181;                              G="TTE" },
    @{ Row=663; A="UIC"; B="081553125";  C=44181;                              G="TTE" },
    @{ Row=664; A="UIC"; B="080715006";  C=44181;                              G="TTE" },
    @{ Row=665; A="UIC"; B="080708793";  C=44181;                              G="TTE" },
    @{ Row=666; A="UIC"; B="081528442";  C=44181;                              G="TTE" },
    @{ Row=667; A="UIC"; B="070931852";  C=44181;                              G="TTE" },
    @{ Row=668; A="UIC"; B="080905668";  C=44181;                              G="TTE" },
    @{ Row=669; A="UIC"; B="080073902";  C=44181; D="AF";     E="ESRD";        G="TEE";    H="bedside" },
    @{ Row=670; A="UIC"; B="080073902";  C=44181; D="AF";     E="ESRD";        G="TEE" },
    @{ Row=671; A="UIC"; B="080073902";  C=44181; D="AF";     E="ESRD";        G="DCCV";   H="200J" },
    @{ Row=672; A="UIC"; B="051582518";  C=44181;                              G="TTE";    I="ASD" },
    @{ Row=673; A="UIC"; B="077069367";  C=44181; D="angina";                  G="stress"; H="treadmill" },
    @{ Row=674; A="UIC"; B="081276913";  C=44181; D="angina";                  G="stress"; H="treadmill" },
    @{ Row=675; A="UIC"; B="080708793";  C=44181; D="angina";                  G="stress"; H="treadmill" },
    @{ Row=676; A="UIC"; B="076821909";  C=44181; D="pre-op";                  G="stress"; H="treadmill" },
    @{ Row=677; A="UIC"; B="081546846";  C=44181; D="angina";                  G="stress"; H="treadmill" },
    @{ Row=678; A="JBVA"; B="W6561";     C=44183; D="CAD";    E="tachycardia";             H="clinic" }
)

foreach ($rec in $data) {
    $r = $rec.Row
    if ($rec.ContainsKey("A")) { $ws.Cells.Item($r, 1).Value = $rec.A }
    if ($rec.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $rec.B }
    if ($rec.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $rec.C }
    if ($rec.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $rec.D }
    if ($rec.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $rec.E }
    if ($rec.ContainsKey("F")) { $ws.Cells.Item($r, 6).Value = $rec.F }
    if ($rec.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $rec.G }
    if ($rec.ContainsKey("H")) { $ws.Cells.Item($r, 8).Value = $rec.H }
    if ($rec.ContainsKey("I")) { $ws.Cells.Item($r, 9).Value = $rec.I }
}

# Make sure the newly written DATE column (C) carries the same date
# number-format / style as the rest of the column, by copying the format
# from an existing date cell rather than letting a brand-new style get
# created.
$ws.Range("C646").Copy() | Out-Null
$ws.Range("C648:C678").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- View state: scroll/freeze to show the newly added rows ---
$ws.Activate()
$ws.Range("E679").Select()

Write-Host "edit complete"
